$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1661.9131
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1661.9131
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4985.7393
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -5321.7393
$ws.Range("H107").Value = 1075.45
$ws.Range("I107").Value = 1115.6154
$ws.Range("K107").Value = 1115.6154
$ws.Range("M107").Value = 804.3846000000001
$ws.Range("H111").Value = 3649.7856
$ws.Range("J111").Value = 5675
$ws.Range("L111").Value = 17025
$ws.Range("N111").Value = -23159
$ws.Range("H131").Value = 948.7
$ws.Range("I131").Value = 893.4211
$ws.Range("J131").Value = 1999
$ws.Range("K131").Value = 2680.2633
$ws.Range("L131").Value = 5997
$ws.Range("M131").Value = 2359.7367
$ws.Range("N131").Value = -16077
$ws.Range("H137").Value = 1933.4147
$ws.Range("I137").Value = 1141.6818
$ws.Range("J137").Value = 2850.158
$ws.Range("K137").Value = 3425.0454
$ws.Range("L137").Value = 8550.474
$ws.Range("M137").Value = -875.0454
$ws.Range("N137").Value = -13650.474
$ws.Range("H138").Value = 2896.7
$ws.Range("I138").Value = 1419.4242
$ws.Range("J138").Value = 3624.3135
$ws.Range("K138").Value = 4258.2726
$ws.Range("L138").Value = 10872.9405
$ws.Range("M138").Value = 881.7273999999998
$ws.Range("N138").Value = -21152.9405

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2380.4
$ws.Range("I63").Value = 1869.7059
$ws.Range("J63").Value = 3465.625
$ws.Range("K63").Value = 1869.7059
$ws.Range("L63").Value = 3465.625
$ws.Range("M63").Value = -1183.7059
$ws.Range("N63").Value = -4837.625
$ws.Range("H66").Value = 2380.4
$ws.Range("I66").Value = 1869.7059
$ws.Range("J66").Value = 3465.625
$ws.Range("K66").Value = 9348.529500000001
$ws.Range("L66").Value = 17328.125
$ws.Range("M66").Value = -5916.529500000001
$ws.Range("N66").Value = -24192.125
$ws.Range("H107").Value = 23942.666
$ws.Range("J107").Value = 23942.666
$ws.Range("L107").Value = 23942.666
$ws.Range("N107").Value = -31622.666
$ws.Range("H110").Value = 1345.5883
$ws.Range("I110").Value = 1678
$ws.Range("K110").Value = 1678
$ws.Range("M110").Value = 367

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H31").Value = 3457.0356
$ws.Range("I31").Value = 1684.8422
$ws.Range("J31").Value = 7198.3335
$ws.Range("K31").Value = 1684.8422
$ws.Range("L31").Value = 7198.3335
$ws.Range("M31").Value = -1389.8422
$ws.Range("N31").Value = -7788.3335
$ws.Range("H34").Value = 3457.0356
$ws.Range("I34").Value = 1684.8422
$ws.Range("J34").Value = 7198.3335
$ws.Range("K34").Value = 1684.8422
$ws.Range("L34").Value = 7198.3335
$ws.Range("M34").Value = -1482.8422
$ws.Range("N34").Value = -7602.3335
$ws.Range("H131").Value = 47390
$ws.Range("J131").Value = 47390
$ws.Range("L131").Value = 47390
$ws.Range("N131").Value = -57470

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 354.82352
$ws.Range("I14").Value = 354.82352
$ws.Range("K14").Value = 1064.47056
$ws.Range("M14").Value = -891.47056
$ws.Range("H64").Value = 8640.571
$ws.Range("I64").Value = 816.6667
$ws.Range("J64").Value = 10774.363
$ws.Range("K64").Value = 2450.0001
$ws.Range("L64").Value = 32323.089
$ws.Range("M64").Value = -2180.0001
$ws.Range("N64").Value = -32863.089
$ws.Range("H67").Value = 8640.571
$ws.Range("I67").Value = 816.6667
$ws.Range("J67").Value = 10774.363
$ws.Range("K67").Value = 2450.0001
$ws.Range("L67").Value = 32323.089
$ws.Range("M67").Value = -1514.0001
$ws.Range("N67").Value = -34195.089
$ws.Range("H94").Value = 9035.333000000001
$ws.Range("I94").Value = 3581
$ws.Range("J94").Value = 19944
$ws.Range("K94").Value = 10743
$ws.Range("L94").Value = 59832
$ws.Range("M94").Value = -10067
$ws.Range("N94").Value = -61184
$ws.Range("H106").Value = 5825.3335
$ws.Range("J106").Value = 5825.3335
$ws.Range("L106").Value = 17476.0005
$ws.Range("N106").Value = -19368.0005
$ws.Range("H107").Value = 843.0303
$ws.Range("I107").Value = 650.2
$ws.Range("J107").Value = 1139.6923
$ws.Range("K107").Value = 1950.6
$ws.Range("L107").Value = 3419.0769
$ws.Range("M107").Value = -30.60000000000014
$ws.Range("N107").Value = -7259.0769
$ws.Range("H109").Value = 2927.8
$ws.Range("I109").Value = 936.1667
$ws.Range("J109").Value = 4255.5557
$ws.Range("K109").Value = 2808.5001
$ws.Range("L109").Value = 12766.6671
$ws.Range("M109").Value = -1768.5001
$ws.Range("N109").Value = -14846.6671
$ws.Range("H112").Value = 46708360
$ws.Range("I112").Value = 900
$ws.Range("J112").Value = 53077560
$ws.Range("K112").Value = 2700
$ws.Range("L112").Value = 159232680
$ws.Range("M112").Value = -1592
$ws.Range("N112").Value = -159234896
$ws.Range("H115").Value = 1420
$ws.Range("I115").Value = 1275
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 3825
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = -2650
$ws.Range("N115").Value = -8350
$ws.Range("H119").Value = 4486.75
$ws.Range("I119").Value = 3171.6365
$ws.Range("J119").Value = 7380
$ws.Range("K119").Value = 9514.9095
$ws.Range("L119").Value = 22140
$ws.Range("M119").Value = -4676.9095
$ws.Range("N119").Value = -31816
$ws.Range("H132").Value = 602813.9399999999
$ws.Range("I132").Value = 1013569.7
$ws.Range("K132").Value = 9122127.299999999
$ws.Range("M132").Value = -9119597.299999999
$ws.Range("H134").Value = 2977.4075
$ws.Range("I134").Value = 2525.7896
$ws.Range("J134").Value = 4050
$ws.Range("K134").Value = 7577.3688
$ws.Range("L134").Value = 12150
$ws.Range("M134").Value = -2507.3688
$ws.Range("N134").Value = -22290

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7222389
$ws.Range("I11").Value = 7500083.5
$ws.Range("J11").Value = 6667000
$ws.Range("K11").Value = 7500083.5
$ws.Range("L11").Value = 6667000
$ws.Range("M11").Value = -7499944.5
$ws.Range("N11").Value = -6667278

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3611.1538
$ws.Range("I61").Value = 2493.125
$ws.Range("K61").Value = 2493.125
$ws.Range("M61").Value = -2291.125
$ws.Range("H113").Value = 3611.1538
$ws.Range("I113").Value = 2493.125
$ws.Range("K113").Value = 2493.125
$ws.Range("M113").Value = -323.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1000839.3
$ws.Range("I100").Value = 457.14285
$ws.Range("J100").Value = 3335064.2
$ws.Range("K100").Value = 914.2857
$ws.Range("L100").Value = 6670128.4
$ws.Range("M100").Value = -373.2857
$ws.Range("N100").Value = -6671210.4
